$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column G: copy formatting from column F header cell, then set values ---
$ws.Cells.Item(1, 6).Copy($ws.Cells.Item(1, 7))

# --- Update header row (B1:G1) with new recipe names ---
$ws.Cells.Item(1, 2).Value = "Recipe Lomo Saltado with beef"
$ws.Cells.Item(1, 3).Value = "Recipe Lomo Saltado with chicken"
$ws.Cells.Item(1, 4).Value = "Recipe Lomo Saltado with chickpeas"
$ws.Cells.Item(1, 5).Value = "Recipe Lomo Saltado with mushrooms"
$ws.Cells.Item(1, 6).Value = "Recipe Lomo Saltado with salmon"
$ws.Cells.Item(1, 7).Value = "Recipe Lomo Saltado with tofu"

# --- Update numeric data for rows 2-15, columns B-G ---
$ws.Cells.Item(2, 2).Value = 13.7190236
$ws.Cells.Item(2, 3).Value = 6.9636386
$ws.Cells.Item(2, 4).Value = 8.8705436
$ws.Cells.Item(2, 5).Value = 4.4898036
$ws.Cells.Item(2, 6).Value = 5.6031886
$ws.Cells.Item(2, 7).Value = 5.513553600000001
$ws.Cells.Item(3, 2).Value = 3.0133508
$ws.Cells.Item(3, 3).Value = 1.5405158
$ws.Cells.Item(3, 4).Value = 2.4227608
$ws.Cells.Item(3, 5).Value = 0.7474508
$ws.Cells.Item(3, 6).Value = 1.0300058
$ws.Cells.Item(3, 7).Value = 1.0177208
$ws.Cells.Item(4, 2).Value = 76.48186
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(5, 2).Value = 8.503787599999999
$ws.Cells.Item(5, 3).Value = 2.2529976
$ws.Cells.Item(5, 4).Value = 1.1123126
$ws.Cells.Item(5, 5).Value = 1.0977526
$ws.Cells.Item(5, 6).Value = 1.3402676
$ws.Cells.Item(5, 7).Value = 1.5331876
$ws.Cells.Item(6, 2).Value = 0.8544899999999999
$ws.Cells.Item(6, 3).Value = 0.807625
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0.372645
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(7, 2).Value = 20.7405214
$ws.Cells.Item(7, 3).Value = 3.1247414
$ws.Cells.Item(7, 4).Value = 3.3144764
$ws.Cells.Item(7, 5).Value = 2.8981514
$ws.Cells.Item(7, 6).Value = 4.456071400000001
$ws.Cells.Item(7, 7).Value = 3.0442064
$ws.Cells.Item(8, 2).Value = 0.6889980000000001
$ws.Cells.Item(8, 3).Value = 0.311803
$ws.Cells.Item(8, 4).Value = 0.111603
$ws.Cells.Item(8, 5).Value = 0.138903
$ws.Cells.Item(8, 6).Value = 0.120703
$ws.Cells.Item(8, 7).Value = 0.4728730000000001
$ws.Cells.Item(9, 2).Value = 0.3578433999999999
$ws.Cells.Item(9, 3).Value = 0.3259933999999999
$ws.Cells.Item(9, 4).Value = 0.2440934
$ws.Cells.Item(9, 5).Value = 0.2750333999999999
$ws.Cells.Item(9, 6).Value = 0.2509184
$ws.Cells.Item(9, 7).Value = 0.2809483999999999
$ws.Cells.Item(10, 2).Value = 0.3899102
$ws.Cells.Item(10, 3).Value = 0.3739852000000001
$ws.Cells.Item(10, 4).Value = 0.4381402
$ws.Cells.Item(10, 5).Value = 0.2961802000000001
$ws.Cells.Item(10, 6).Value = 0.3052802000000001
$ws.Cells.Item(10, 7).Value = 0.3580602
$ws.Cells.Item(11, 2).Value = 0.1350142
$ws.Cells.Item(11, 3).Value = 0.1409292
$ws.Cells.Item(11, 4).Value = 0.07813920000000001
$ws.Cells.Item(11, 5).Value = 0.06721920000000001
$ws.Cells.Item(11, 6).Value = 0.0785942
$ws.Cells.Item(11, 7).Value = 0.1832442
$ws.Cells.Item(12, 2).Value = 0.107303026
$ws.Cells.Item(12, 3).Value = 0.050473526
$ws.Cells.Item(12, 4).Value = 0.029880226
$ws.Cells.Item(12, 5).Value = 0.023592126
$ws.Cells.Item(12, 6).Value = 0.032551076
$ws.Cells.Item(12, 7).Value = 0.024019826
$ws.Cells.Item(13, 2).Value = 0.103184908
$ws.Cells.Item(13, 3).Value = 0.035890408
$ws.Cells.Item(13, 4).Value = 0.029265608
$ws.Cells.Item(13, 5).Value = 0.023282358
$ws.Cells.Item(13, 6).Value = 0.06432790799999999
$ws.Cells.Item(13, 7).Value = 0.024952208
$ws.Cells.Item(14, 2).Value = 1432.1124
$ws.Cells.Item(14, 3).Value = 1218.7174
$ws.Cells.Item(14, 4).Value = 1201.4274
$ws.Cells.Item(14, 5).Value = 1073.5724
$ws.Cells.Item(14, 6).Value = 1634.1324
$ws.Cells.Item(14, 7).Value = 1099.0524
$ws.Cells.Item(15, 2).Value = 40831.37280000001
$ws.Cells.Item(15, 3).Value = 35218.94780000002
$ws.Cells.Item(15, 4).Value = 39697.05780000002
$ws.Cells.Item(15, 5).Value = 33003.0978
$ws.Cells.Item(15, 6).Value = 37431.61280000001
$ws.Cells.Item(15, 7).Value = 33379.38280000001
